$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A45").Value = "Offerte voor zakelijke samenwerking"
$ws.Range("B45").Value = "mailmind.test@zohomail.eu"
$ws.Range("C45").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$ws.Range("D45").Value = "Offerte / Prijsaanvraag"
$ws.Range("F45").Value = "2025-06-19 22:32:11"
$ws.Range("G45").Value = "Nee"

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("B5").Value = 5
$dash.Range("A6").Value = "Klacht / Probleem"
$dash.Range("B6").Value = 4

$dCond = $ws.Range("D2:D44").FormatConditions.Item(1)
$dCond.ModifyAppliesToRange($ws.Range("D2:D45"))

$gCond = $ws.Range("G2:G44").FormatConditions.Item(1)
$gCond.ModifyAppliesToRange($ws.Range("G2:G45"))
